$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared SHAP split-scores string (column Z) for all data rows
$ws.Range("Z2:Z17").Value = "[0.63492063 0.55357143 0.68415638 0.70884774 0.68979592]"

# Row 2
$ws.Range("B2").Value = 0.005730915069580078
$ws.Range("C2").Value = 0.003000456279259475
$ws.Range("D2").Value = 0.001888942718505859
$ws.Range("E2").Value = 0.0004769262564387624
$ws.Range("R2").Value = 0.5964912280701754
$ws.Range("S2").Value = 0.4
$ws.Range("T2").Value = 0.5862068965517241
$ws.Range("U2").Value = 0.6296296296296297
$ws.Range("V2").Value = 0.6071428571428571
$ws.Range("W2").Value = 0.5638941222788774
$ws.Range("X2").Value = 0.08320255346235392
$ws.Range("Y2").Value = 9
$ws.Range("AA2").Value = 0.6542584194171497
$ws.Range("AB2").Value = 0.0559439279691001

# Row 3
$ws.Range("B3").Value = 0.004150629043579102
$ws.Range("C3").Value = 0.001408240015526683
$ws.Range("D3").Value = 0.001403999328613281
$ws.Range("E3").Value = 0.0002202478469690135
$ws.Range("R3").Value = 0.5483870967741935
$ws.Range("S3").Value = 0.4745762711864406
$ws.Range("T3").Value = 0.6567164179104479
$ws.Range("U3").Value = 0.5555555555555556
$ws.Range("V3").Value = 0.6567164179104477
$ws.Range("W3").Value = 0.5783903518674171
$ws.Range("X3").Value = 0.06995560282062323
$ws.Range("Y3").Value = 3
$ws.Range("AA3").Value = 0.6542584194171497
$ws.Range("AB3").Value = 0.0559439279691001

# Row 4
$ws.Range("B4").Value = 0.003569650650024414
$ws.Range("C4").Value = 0.001197459058299865
$ws.Range("D4").Value = 0.001888656616210938
$ws.Range("E4").Value = 0.0006248975173018375
$ws.Range("R4").Value = 0.5818181818181818
$ws.Range("S4").Value = 0.4081632653061225
$ws.Range("T4").Value = 0.5714285714285714
$ws.Range("U4").Value = 0.6153846153846153
$ws.Range("V4").Value = 0.5660377358490566
$ws.Range("W4").Value = 0.5485664739573095
$ws.Range("X4").Value = 0.0722652358909711
$ws.Range("Y4").Value = 15
$ws.Range("AA4").Value = 0.6542584194171497
$ws.Range("AB4").Value = 0.0559439279691001

# Row 5
$ws.Range("B5").Value = 0.004147768020629883
$ws.Range("C5").Value = 0.0002988488396109322
$ws.Range("D5").Value = 0.002070426940917969
$ws.Range("E5").Value = 0.0005159780835631233
$ws.Range("R5").Value = 0.5333333333333333
$ws.Range("S5").Value = 0.4285714285714285
$ws.Range("T5").Value = 0.6153846153846153
$ws.Range("U5").Value = 0.5185185185185185
$ws.Range("V5").Value = 0.6562499999999999
$ws.Range("W5").Value = 0.5504115791615791
$ws.Range("X5").Value = 0.0794696455952751
$ws.Range("Y5").Value = 11
$ws.Range("AA5").Value = 0.6542584194171497
$ws.Range("AB5").Value = 0.0559439279691001

# Row 6
$ws.Range("B6").Value = 0.007025861740112304
$ws.Range("C6").Value = 0.003903439402866177
$ws.Range("D6").Value = 0.001665067672729492
$ws.Range("E6").Value = 0.001033513090673748
$ws.Range("R6").Value = 0.5762711864406779
$ws.Range("S6").Value = 0.4
$ws.Range("T6").Value = 0.6461538461538462
$ws.Range("U6").Value = 0.6545454545454545
$ws.Range("V6").Value = 0.6774193548387097
$ws.Range("W6").Value = 0.5908779683957377
$ws.Range("X6").Value = 0.1012384708695192
$ws.Range("Y6").Value = 1
$ws.Range("AA6").Value = 0.6542584194171497
$ws.Range("AB6").Value = 0.0559439279691001

# Row 7
$ws.Range("B7").Value = 0.009363079071044922
$ws.Range("C7").Value = 0.006019449554514123
$ws.Range("D7").Value = 0.003815174102783203
$ws.Range("E7").Value = 0.002137701935921446
$ws.Range("R7").Value = 0.5483870967741935
$ws.Range("S7").Value = 0.4745762711864406
$ws.Range("T7").Value = 0.6567164179104479
$ws.Range("U7").Value = 0.5555555555555556
$ws.Range("V7").Value = 0.6567164179104477
$ws.Range("W7").Value = 0.5783903518674171
$ws.Range("X7").Value = 0.06995560282062323
$ws.Range("Y7").Value = 3
$ws.Range("AA7").Value = 0.6542584194171497
$ws.Range("AB7").Value = 0.0559439279691001

# Row 8
$ws.Range("B8").Value = 0.006037092208862305
$ws.Range("C8").Value = 0.002059634082720732
$ws.Range("D8").Value = 0.00671391487121582
$ws.Range("E8").Value = 0.004049177161724995
$ws.Range("R8").Value = 0.5517241379310344
$ws.Range("S8").Value = 0.3703703703703704
$ws.Range("T8").Value = 0.6333333333333334
$ws.Range("U8").Value = 0.6181818181818182
$ws.Range("V8").Value = 0.6440677966101694
$ws.Range("W8").Value = 0.5635354912853452
$ws.Range("X8").Value = 0.1017798447679304
$ws.Range("Y8").Value = 10
$ws.Range("AA8").Value = 0.6542584194171497
$ws.Range("AB8").Value = 0.0559439279691001

# Row 9
$ws.Range("B9").Value = 0.01439299583435059
$ws.Range("C9").Value = 0.007901324228454326
$ws.Range("D9").Value = 0.003690242767333984
$ws.Range("E9").Value = 0.003493173771322045
$ws.Range("R9").Value = 0.5333333333333333
$ws.Range("S9").Value = 0.4285714285714285
$ws.Range("T9").Value = 0.6153846153846153
$ws.Range("U9").Value = 0.5185185185185185
$ws.Range("V9").Value = 0.6562499999999999
$ws.Range("W9").Value = 0.5504115791615791
$ws.Range("X9").Value = 0.0794696455952751
$ws.Range("Y9").Value = 11
$ws.Range("AA9").Value = 0.6542584194171497
$ws.Range("AB9").Value = 0.0559439279691001

# Row 10
$ws.Range("B10").Value = 0.01257500648498535
$ws.Range("C10").Value = 0.005673234197555354
$ws.Range("D10").Value = 0.005941152572631836
$ws.Range("E10").Value = 0.001886994350153504
$ws.Range("R10").Value = 0.603174603174603
$ws.Range("S10").Value = 0.4285714285714285
$ws.Range("T10").Value = 0.6197183098591549
$ws.Range("U10").Value = 0.5818181818181818
$ws.Range("V10").Value = 0.6666666666666667
$ws.Range("W10").Value = 0.5799898380180071
$ws.Range("X10").Value = 0.08069651169061996
$ws.Range("Y10").Value = 2
$ws.Range("AA10").Value = 0.6542584194171497
$ws.Range("AB10").Value = 0.0559439279691001

# Row 11
$ws.Range("B11").Value = 0.01639413833618164
$ws.Range("C11").Value = 0.005223813367791615
$ws.Range("D11").Value = 0.004514837265014648
$ws.Range("E11").Value = 0.003328428279075724
$ws.Range("R11").Value = 0.5483870967741935
$ws.Range("S11").Value = 0.4745762711864406
$ws.Range("T11").Value = 0.6567164179104479
$ws.Range("U11").Value = 0.5555555555555556
$ws.Range("V11").Value = 0.6567164179104477
$ws.Range("W11").Value = 0.5783903518674171
$ws.Range("X11").Value = 0.06995560282062323
$ws.Range("Y11").Value = 3
$ws.Range("AA11").Value = 0.6542584194171497
$ws.Range("AB11").Value = 0.0559439279691001

# Row 12
$ws.Range("B12").Value = 0.02494139671325684
$ws.Range("C12").Value = 0.0197527262375071
$ws.Range("D12").Value = 0.007191038131713868
$ws.Range("E12").Value = 0.006097777330149007
$ws.Range("R12").Value = 0.5806451612903226
$ws.Range("S12").Value = 0.4
$ws.Range("T12").Value = 0.6031746031746033
$ws.Range("U12").Value = 0.5818181818181818
$ws.Range("V12").Value = 0.6666666666666667
$ws.Range("W12").Value = 0.5664609225899548
$ws.Range("X12").Value = 0.08892172094156747
$ws.Range("Y12").Value = 8
$ws.Range("AA12").Value = 0.6542584194171497
$ws.Range("AB12").Value = 0.0559439279691001

# Row 13
$ws.Range("B13").Value = 0.009238243103027344
$ws.Range("C13").Value = 0.005730825747515498
$ws.Range("D13").Value = 0.0082305908203125
$ws.Range("E13").Value = 0.01073906746615409
$ws.Range("R13").Value = 0.5333333333333333
$ws.Range("S13").Value = 0.4285714285714285
$ws.Range("T13").Value = 0.6153846153846153
$ws.Range("U13").Value = 0.5185185185185185
$ws.Range("V13").Value = 0.6562499999999999
$ws.Range("W13").Value = 0.5504115791615791
$ws.Range("X13").Value = 0.0794696455952751
$ws.Range("Y13").Value = 11
$ws.Range("AA13").Value = 0.6542584194171497
$ws.Range("AB13").Value = 0.0559439279691001

# Row 14
$ws.Range("B14").Value = 0.01393752098083496
$ws.Range("C14").Value = 0.01233294558399919
$ws.Range("D14").Value = 0.002502298355102539
$ws.Range("E14").Value = 0.001558685552181106
$ws.Range("R14").Value = 0.5483870967741935
$ws.Range("S14").Value = 0.4745762711864406
$ws.Range("T14").Value = 0.6567164179104479
$ws.Range("U14").Value = 0.5454545454545454
$ws.Range("V14").Value = 0.6567164179104477
$ws.Range("W14").Value = 0.576370149847215
$ws.Range("X14").Value = 0.07072745601311242
$ws.Range("Y14").Value = 7
$ws.Range("AA14").Value = 0.6542584194171497
$ws.Range("AB14").Value = 0.0559439279691001

# Row 15
$ws.Range("B15").Value = 0.0150115966796875
$ws.Range("C15").Value = 0.009538176237436949
$ws.Range("D15").Value = 0.008159112930297852
$ws.Range("E15").Value = 0.00798967269282398
$ws.Range("R15").Value = 0.5483870967741935
$ws.Range("S15").Value = 0.4745762711864406
$ws.Range("T15").Value = 0.6567164179104479
$ws.Range("U15").Value = 0.5555555555555556
$ws.Range("V15").Value = 0.6567164179104477
$ws.Range("W15").Value = 0.5783903518674171
$ws.Range("X15").Value = 0.06995560282062323
$ws.Range("Y15").Value = 3
$ws.Range("AA15").Value = 0.6542584194171497
$ws.Range("AB15").Value = 0.0559439279691001

# Row 16
$ws.Range("B16").Value = 0.003835773468017578
$ws.Range("C16").Value = 0.002085516593103265
$ws.Range("D16").Value = 0.001375389099121094
$ws.Range("E16").Value = 0.0006762739428131437
$ws.Range("R16").Value = 0.5333333333333333
$ws.Range("S16").Value = 0.4285714285714285
$ws.Range("T16").Value = 0.6153846153846153
$ws.Range("U16").Value = 0.5185185185185185
$ws.Range("V16").Value = 0.6363636363636364
$ws.Range("W16").Value = 0.5464343064343063
$ws.Range("X16").Value = 0.07440969584019633
$ws.Range("Y16").Value = 16
$ws.Range("AA16").Value = 0.6542584194171497
$ws.Range("AB16").Value = 0.0559439279691001

# Row 17
$ws.Range("B17").Value = 0.002525615692138672
$ws.Range("C17").Value = 0.001056815898788723
$ws.Range("D17").Value = 0.001015615463256836
$ws.Range("E17").Value = 0.0003878686888305205
$ws.Range("R17").Value = 0.5333333333333333
$ws.Range("S17").Value = 0.4285714285714285
$ws.Range("T17").Value = 0.6153846153846153
$ws.Range("U17").Value = 0.5185185185185185
$ws.Range("V17").Value = 0.6562499999999999
$ws.Range("W17").Value = 0.5504115791615791
$ws.Range("X17").Value = 0.0794696455952751
$ws.Range("Y17").Value = 11
$ws.Range("AA17").Value = 0.6542584194171497
$ws.Range("AB17").Value = 0.0559439279691001
